$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Apply-DataStyle {
    param($addr, $center, $border, $dateFmt)
    $r = $ws.Range($addr)
    $r.Font.Name = "Aptos Narrow"
    $r.Font.Color = 0
    if ($border -eq 1) {
        $r.Borders.LineStyle = 1
    }
    if ($center -eq 1) {
        $r.HorizontalAlignment = -4108
    }
    if ($dateFmt -eq 1) {
        $r.NumberFormat = "m/d/yyyy"
    }
}

# ---- New rows 105-111 (daily entries for 1/29 - 1/30 plus totals) ----

# Row 105: date + Domm + 0.25
$ws.Range("A105").Value = 45687
$ws.Range("B105").Value = "Domm"
$ws.Range("D105").Value = 0.25

# Row 106: Meeting / Reconsile / 1
$ws.Range("B106").Value = "Meeting"
$ws.Range("C106").Value = "Reconsile"
$ws.Range("D106").Value = 1

# Row 107: General Discussion / 0.25
$ws.Range("C107").Value = "General Discussion"
$ws.Range("D107").Value = 0.25

# Row 108: Study / Documentation / 2.5
$ws.Range("B108").Value = "Study"
$ws.Range("C108").Value = "Documentation"
$ws.Range("D108").Value = 2.5

# Row 109: Reconcile Corrections / 1
$ws.Range("C109").Value = "Reconcile Corrections"
$ws.Range("D109").Value = 1

# Row 110: Reconcile Revision / 3
$ws.Range("C110").Value = "Reconcile Revision"
$ws.Range("D110").Value = 3

# Row 111: Total label + SUM formula
$ws.Range("B111").Value = "Total"
$ws.Range("D111").Formula = "=SUM(D104:D110)"

# ---- Formatting to match the rest of the table (best effort) ----

Apply-DataStyle "A105" 0 1 1
Apply-DataStyle "B105" 0 1 0
Apply-DataStyle "C105" 0 1 0
Apply-DataStyle "D105" 1 1 0

Apply-DataStyle "A106" 0 1 1
Apply-DataStyle "B106" 0 1 0
Apply-DataStyle "C106" 0 1 0
Apply-DataStyle "D106" 1 1 0

Apply-DataStyle "A107" 0 1 0
Apply-DataStyle "B107" 0 1 0
Apply-DataStyle "C107" 0 1 0
Apply-DataStyle "D107" 1 1 0

Apply-DataStyle "A108" 0 1 0
Apply-DataStyle "B108" 0 1 0
Apply-DataStyle "C108" 0 1 0
Apply-DataStyle "D108" 1 1 0

Apply-DataStyle "A109" 0 1 0
Apply-DataStyle "B109" 0 1 0
Apply-DataStyle "C109" 0 1 0
Apply-DataStyle "D109" 1 1 0

Apply-DataStyle "A110" 0 1 0
Apply-DataStyle "B110" 0 1 0
Apply-DataStyle "C110" 0 1 0
Apply-DataStyle "D110" 1 1 0

Apply-DataStyle "A111" 0 1 0
Apply-DataStyle "B111" 0 1 0
Apply-DataStyle "C111" 0 1 0
Apply-DataStyle "D111" 1 1 0

# ---- Sheet view / selection update ----
$ws.Range("A104:D112").Select()
$wb.Windows.Item(1).ScrollRow = 87
